# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.585.39"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.814.32"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'229.03"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("E6").Value = "  +4.53%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'34.98"
$ws.Range("E8").Value = "  +7.16%  "
$ws.Range("E9").Value = "  +1.68%  "
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "2.076.92"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.873.41"
$ws.Range("E13").Value = "  +4.28%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'11.24"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "'0.649"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'4.47"
$ws.Range("E16").Value = "  +3.75%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "34.560.43"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "'69.29"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "'246.19"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D24").Value = "'174.02"
$ws.Range("E24").Value = "  +3.35%  "
$ws.Range("D25").Value = "'2.13"
$ws.Range("E25").Value = "  +2.50%  "
$ws.Range("D26").Value = "'7.96"
$ws.Range("E26").Value = "  +9.14%  "
$ws.Range("D27").Value = "'16.87"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "'4.03"
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("D31").Value = "'0.0534"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("D32").Value = "'3.87"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("D36").Value = "1.398.86"
$ws.Range("E36").Value = "  -2.27%  "
$ws.Range("D37").Value = "'2.50"
$ws.Range("E37").Value = "  -3.85%  "
$ws.Range("D38").Value = "'1.06"
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "'83.91"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("D43").Value = "'2.39"
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("D44").Value = "'13.56"
$ws.Range("E44").Value = "  -2.33%  "
$ws.Range("E45").Value = "  +3.99%  "
$ws.Range("D46").Value = "'0.0514"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").Value = "1.976.43"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("D49").Value = "'105.09"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("D50").Value = "0.0₆0131"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("E51").Value = "  -0.04%  "
